$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.372.67"
$ws.Range("E2").Value = "  +9.19%  "
$ws.Range("D3").Value = "3.162.10"
$ws.Range("E3").Value = "  +6.54%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'591.72"
$ws.Range("E5").Value = "  +5.22%  "
$ws.Range("D6").Value = "'148.57"
$ws.Range("E6").Value = "  +9.10%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.157.53"
$ws.Range("E8").Value = "  +6.55%  "
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  +20.80%  "
$ws.Range("D11").Value = "'5.81"
$ws.Range("E11").Value = "  +10.69%  "
$ws.Range("E12").Value = "  +5.89%  "
$ws.Range("E13").Value = "  +11.83%  "
$ws.Range("D14").Value = "'36.12"
$ws.Range("E14").Value = "  +8.08%  "
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "3.684.44"
$ws.Range("E16").Value = "  +6.54%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "64.276.72"
$ws.Range("E17").Value = "  +8.99%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.20"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").Value = "3.159.48"
$ws.Range("E19").Value = "  +6.43%  "
$ws.Range("D20").Value = "'478.27"
$ws.Range("E20").Value = "  +10.45%  "
$ws.Range("D21").Value = "'14.33"
$ws.Range("E21").Value = "  +5.90%  "
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("E23").Value = "  +10.28%  "
$ws.Range("D24").Value = "'13.39"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "'82.90"
$ws.Range("E25").Value = "  +4.14%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'8.72"
$ws.Range("E27").Value = "  +13.91%  "
$ws.Range("D28").Value = "'2.25"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("E29").Value = "  +7.11%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "'6.93"
$ws.Range("E31").Value = "  +13.00%  "
$ws.Range("D32").Value = "'27.40"
$ws.Range("E32").Value = "  +7.40%  "
$ws.Range("D33").Value = "'0.110"
$ws.Range("E33").Value = "  +7.01%  "
$ws.Range("D34").Value = "0.0₃0889"
$ws.Range("E34").Value = "  +17.64%  "
$ws.Range("E35").Value = "  +19.63%  "
$ws.Range("E36").Value = "  +8.31%  "
$ws.Range("D37").Value = "'3.43"
$ws.Range("E37").Value = "  +25.35%  "
$ws.Range("D38").Value = "'6.19"
$ws.Range("E38").Value = "  +5.87%  "
$ws.Range("D39").Value = "'50.93"
$ws.Range("E39").Value = "  +5.20%  "
$ws.Range("D40").Value = "'454.22"
$ws.Range("E40").Value = "  +15.06%  "
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0376"
$ws.Range("E42").Value = "  +7.66%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.948.55"
$ws.Range("E43").Value = "  +8.59%  "
$ws.Range("E44").Value = "  +14.70%  "
$ws.Range("E45").Value = "  +6.91%  "
$ws.Range("E46").Value = "  +13.72%  "
$ws.Range("D47").Value = "'35.56"
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("D49").Value = "'123.60"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").Value = "'25.13"
$ws.Range("E51").Value = "  +9.35%  "
